$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column C stays as text so the date-like strings are not
# reinterpreted as date serial numbers.
$ws.Range("C2:C9").NumberFormat = "@"

# Row 2 - Pernambuco
$ws.Range("C2").Value = "01/04/2025"
$ws.Range("D2").Value = 10.4

# Row 3 - Bahia
$ws.Range("C3").Value = "01/04/2025"
$ws.Range("D3").Value = 9.1

# Row 4 - Distrito Federal
$ws.Range("C4").Value = "01/04/2025"
$ws.Range("D4").Value = 8.699999999999999

# Row 5 - Amapá -> Piauí
$ws.Range("A5").Value = "Piauí"
$ws.Range("C5").Value = "01/04/2025"
$ws.Range("D5").Value = 8.5

# Row 6 - Rio Grande do Norte -> Rio de Janeiro
$ws.Range("A6").Value = "Rio de Janeiro"
$ws.Range("C6").Value = "01/04/2025"
$ws.Range("D6").Value = 8.1

# Row 7 - Sergipe
$ws.Range("C7").Value = "01/04/2025"
$ws.Range("D7").Value = 8.1
$ws.Range("E7").Value = "5º"

# Row 8 - Brasil
$ws.Range("C8").Value = "01/04/2025"
$ws.Range("D8").Value = 5.8

# Row 9 - Nordeste
$ws.Range("C9").Value = "01/04/2025"
$ws.Range("D9").Value = 8.199999999999999
